# Add mAP50 / mAP50-95 columns (G, H) with their data, and update
# Recall/Precision (E/F) values that were recalculated after including
# the euclidean distance for tracking trajectories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Headers (row 5)
# ---------------------------------------------------------------
# G5 uses the same header style as the rest of the header row.
$ws.Range("F5").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Value = "mAP50"

# H5 uses a bold+border header style too, but created as a distinct
# (new) style entry -- replicate by copying the header format and then
# re-asserting the fill explicitly (forces a new, equivalent xf record).
$ws.Range("F5").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Interior.ColorIndex = -4142
$ws.Range("H5").Value = "mAP50-95"

# ---------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------

function Set-Val($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 6 (Nano, 640)
Set-Val "E6" 0.928
Set-Val "F6" 0.949
Set-Val "G6" 0.964
Set-Val "H6" 0.701

# Row 7 (Nano, 1200)
Set-Val "E7" 0.881
Set-Val "F7" 0.937
Set-Val "G7" 0.915
Set-Val "H7" 0.701

# Row 8 (Nano, 1920)
Set-Val "E8" 0.931
Set-Val "F8" 0.945
Set-Val "G8" 0.983
Set-Val "H8" 0.729

# Row 9 (Small, 640) -- E9 unchanged
Set-Val "F9" 0.933
Set-Val "G9" 0.961
Set-Val "H9" 0.693

# Row 10 (Small, 1200) -- E10/F10 unchanged
Set-Val "G10" 0.936
Set-Val "H10" 0.682

# Row 11 (Small, 1920, orange fill) -- E11/F11 unchanged
$ws.Range("F11").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
Set-Val "G11" 0.978
$ws.Range("F11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
Set-Val "H11" 0.724

# Row 12 (Medium, 640)
Set-Val "E12" 0.918
Set-Val "F12" 0.922
Set-Val "G12" 0.954
Set-Val "H12" 0.671

# Row 13 (Medium, 1200, yellow fill)
$ws.Range("F13").Copy() | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null
Set-Val "G13" 0.963
$ws.Range("F13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
Set-Val "H13" 0.706
Set-Val "E13" 0.961

# Row 14 (Medium, 1920)
Set-Val "E14" 0.941
Set-Val "G14" 0.97
Set-Val "H14" 0.716

# Row 15 (Large, 640)
Set-Val "E15" 0.866
Set-Val "F15" 0.961
Set-Val "G15" 0.947
Set-Val "H15" 0.704

# Row 16 (Large, 1200)
Set-Val "E16" 0.915
Set-Val "F16" 0.943
Set-Val "G16" 0.954
Set-Val "H16" 0.712

# Row 17 (Large, 1920) -- E17/F17 unchanged
Set-Val "G17" 0.979
Set-Val "H17" 0.737

# Row 18 (ExtraLarge, 640) -- E18/F18 unchanged
Set-Val "G18" 0.947
Set-Val "H18" 0.704

# Row 19 (ExtraLarge, 1200) -- E19 unchanged
Set-Val "F19" 0.925
Set-Val "G19" 0.961
Set-Val "H19" 0.714

# Row 20 (ExtraLarge, 1920, bottom border row) -- E20/F20 unchanged
$ws.Range("C20").Copy() | Out-Null
$ws.Range("G20").PasteSpecial(-4122) | Out-Null
$ws.Range("G20").Interior.ColorIndex = -4142
Set-Val "G20" 0.97
$ws.Range("C20").Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4122) | Out-Null
$ws.Range("H20").Interior.ColorIndex = -4142
Set-Val "H20" 0.72

# Row 22 (bottom styled blank row) -- extend the existing numeric-format
# blank style into the new G column.
$ws.Range("F22").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null

$ws.Range("N16").Select() | Out-Null
